# Apply the "Optuna Attempt (go back with original)" edits to the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("L2").Value = 1.08

$ws1.Range("L3").Value = 0.84

$ws1.Range("D4").Value = 48
$ws1.Range("H4").Value = 2.93
$ws1.Range("L4").Value = 0.9399999999999999

$ws1.Range("D5").Value = 32
$ws1.Range("H5").Value = 2.82

$ws1.Range("D6").Value = 32
$ws1.Range("H6").Value = 1.82
$ws1.Range("I6").Value = "Low"
$ws1.Range("J6").Value = "Normal"
$ws1.Range("L6").Value = 0.8100000000000001

$ws1.Range("D7").Value = 32
$ws1.Range("H7").Value = 0.84
$ws1.Range("I7").Value = "Low"
$ws1.Range("L7").Value = 1.07

$ws1.Range("D8").Value = 32
$ws1.Range("L8").Value = 0.92

$ws1.Range("D9").Value = 32
$ws1.Range("L9").Value = 1

$ws1.Range("D10").Value = 32
$ws1.Range("L10").Value = 0.99

$ws1.Range("D11").Value = 31
$ws1.Range("L11").Value = 0.9

$ws1.Range("D12").Value = 32
$ws1.Range("L12").Value = 1.09

$ws1.Range("D13").Value = 32
$ws1.Range("L13").Value = 0.9399999999999999

$ws1.Range("D14").Value = 31
$ws1.Range("L14").Value = 1.11

$ws1.Range("D15").Value = 30
$ws1.Range("L15").Value = 1.03

$ws1.Range("D16").Value = 31
$ws1.Range("L16").Value = 1

$ws1.Range("D17").Value = 30
$ws1.Range("L17").Value = 0.9399999999999999

# --- Sheet 2: "Summary" ---
# These B-column values are stored as TEXT (not numbers) in the workbook, even
# though they look numeric, so force text formatting before assigning the new
# value (otherwise Excel auto-converts the numeric-looking string to a real
# number). Reset the style back to "Normal" afterwards so we don't leave a
# stray number-format behind on the cell.
$ws2 = $wb.Worksheets.Item("Summary")

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws2.Range("B9") "585"
Set-TextValue $ws2.Range("B10") "333"
Set-TextValue $ws2.Range("B11") "203"
Set-TextValue $ws2.Range("B14") "30"
